$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: add new "x" marks in L2 and M2
$ws.Range("L2").Value = "x"
$ws.Range("M2").Value = "x"

# Row 8 (intern_.all_args_here_fill.R): clear the long comment in B8,
# set row height, and fill C8, E8:M8 (skip D8) with "x"
$ws.Range("B8").Clear()
$ws.Range("C8").Value = "x"
$ws.Range("E8").Value = "x"
$ws.Range("F8").Value = "x"
$ws.Range("G8").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("I8").Value = "x"
$ws.Range("J8").Value = "x"
$ws.Range("K8").Value = "x"
$ws.Range("L8").Value = "x"
$ws.Range("M8").Value = "x"
$ws.Rows.Item(8).RowHeight = 17.25

# Row 13 (intern_.functions_detect.R): fill C13, E13:M13 (skip D13) with "x"
$ws.Range("C13").Value = "x"
$ws.Range("E13").Value = "x"
$ws.Range("F13").Value = "x"
$ws.Range("G13").Value = "x"
$ws.Range("H13").Value = "x"
$ws.Range("I13").Value = "x"
$ws.Range("J13").Value = "x"
$ws.Range("K13").Value = "x"
$ws.Range("L13").Value = "x"
$ws.Range("M13").Value = "x"

# C17: remove the redundant fill-style variant by re-applying the plain
# centered/wrapped style (same visual style index 2)
$ws.Range("C17").HorizontalAlignment = -4108
$ws.Range("C17").VerticalAlignment = -4108
$ws.Range("C17").WrapText = $true

# Update selection to N8
$ws.Range("N8").Select()

$wb.Save()
